$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts old E:O -> F:P) to make room
# for the new "DE PASCALE" list.
$ws.Range("E1").EntireColumn.Insert()

# --- Header row (row 1) ---
# Note: shared-string table order matters for round-trip fidelity, so add
# "DE PASCALE" before "ALLEANZA VERDI SINISTRA".
$ws.Range("E1").Value = "DE PASCALE"
$ws.Range("C1").Value = "ALLEANZA VERDI SINISTRA"

# --- Content fix-ups that are not simple column shifts ---

# Row 6 (ALTERNATIVA POPOLARE): drop the old C value (old "ALLEANZA VERDI E
# SINISTRA" mapping no longer applies); keep only the self-mapping cell,
# which already shifted from M to N.
$ws.Range("C6").ClearContents()

# Row 12 (EUROPA VERDE): value stays mapped to (the renamed) column C, so
# put it back after the automatic shift moved it to D.
$ws.Range("C12").Value = 1
$ws.Range("D12").ClearContents()

# Row 20 (PARTITO PIRATA): old mapping to "ALLEANZA VERDI E SINISTRA" (old
# C, unaffected by the column insert) is now remapped to "astensione"
# (column P) instead.
$ws.Range("C20").ClearContents()
$ws.Range("P20").Value = 1

# Row 28 (BONACCINI PRESIDENTE): remapped from PARTITO DEMOCRATICO (D) to
# DE PASCALE (E) - the new civic-list analogue.
$ws.Range("D28").ClearContents()
$ws.Range("E28").Value = 1

# Row 38 (ITALIA EUROPA INSIEME): the C value changes from 0.33 to 0.3, and
# stays mapped to (the renamed) column C like row 12/54.
$ws.Range("C38").Value = 0.3
$ws.Range("D38").Value = 0.33

# Row 54 (ALLEANZA VERDI E SINISTRA): stays mapped to (the renamed) column C.
$ws.Range("C54").Value = 1
$ws.Range("D54").ClearContents()

# --- Sheet view / layout touch-ups ---
$ws.Range("C2").Select()

$ws.Columns("C:O").ColumnWidth = 4.85546875
$ws.Columns("P:P").ColumnWidth = 5.28515625

$ws.Range("A1:P1").RowHeight = 129
